$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 236, shifting the existing rows
# 236-260 down to 238-262 (same as Excel's native "Insert Rows").
$ws.Rows("236:237").Insert()

# New row 236
$ws.Cells.Item(236,1).Value2 = 3
$ws.Cells.Item(236,2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(236,3).Value2 = "Coquimbo"
$ws.Cells.Item(236,4).Value2 = 44449
$ws.Cells.Item(236,5).Value2 = 5
$ws.Cells.Item(236,6).Value2 = 100112021
$ws.Cells.Item(236,7).Value2 = "Ají"
$ws.Cells.Item(236,8).Value2 = "Americana (o)"
$ws.Cells.Item(236,9).Value2 = "Primera"
$ws.Cells.Item(236,10).Value2 = 30
$ws.Cells.Item(236,11).Value2 = 71000
$ws.Cells.Item(236,12).Value2 = 71000
$ws.Cells.Item(236,13).Value2 = 71000
$ws.Cells.Item(236,14).Value2 = "`$/caja 25 kilos"
$ws.Cells.Item(236,15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(236,16).Value2 = 2840
$ws.Cells.Item(236,17).Value2 = 25
$ws.Cells.Item(236,18).Value2 = "Hortaliza"

# New row 237
$ws.Cells.Item(237,1).Value2 = 3
$ws.Cells.Item(237,2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(237,3).Value2 = "Coquimbo"
$ws.Cells.Item(237,4).Value2 = 44449
$ws.Cells.Item(237,5).Value2 = 5
$ws.Cells.Item(237,6).Value2 = 100112021
$ws.Cells.Item(237,7).Value2 = "Ají"
$ws.Cells.Item(237,8).Value2 = "Inferno"
$ws.Cells.Item(237,9).Value2 = "Primera"
$ws.Cells.Item(237,10).Value2 = 74
$ws.Cells.Item(237,11).Value2 = 41000
$ws.Cells.Item(237,12).Value2 = 42000
$ws.Cells.Item(237,13).Value2 = 41486
$ws.Cells.Item(237,14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(237,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(237,16).Value2 = 2766
$ws.Cells.Item(237,17).Value2 = 15
$ws.Cells.Item(237,18).Value2 = "Hortaliza"
